# 488-RBI-EPP-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-EarlyRePayment-Newcreateloan.xlsx
# "Loan RBI, Variable Instalments"
#
# The "Repayment Schedule" sheet gets a new (blank) column inserted right
# before the existing "Late" column (column N), pushing "Late" from N -> O
# and the already-blank spacer + "Outstanding" column from O/P -> P/Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new blank column at N; this shifts N (Late) -> O and P
# (Outstanding) -> Q, and also takes care of updating the sheet dimension
# and row/col spans to match.
$ws.Columns("N:N").Insert()

# The inserted column inherits the width of the column to its left (M -
# "In Advance") instead of keeping the old bestFit flag/width that used to
# live at N.
$ws.Columns("N:N").ColumnWidth = 10.33

# Match the author's cursor position recorded in the saved view state.
[void]$ws.Range("S8").Select()
